# "cierre de cambios qmas 3 feb"
# - Row 2 (Caren Lorena Gallego Peña / 1000983509) becomes Stefany Jimenez / 1000588257
# - Row 3 (Michell Mesa / 1001315941) is removed entirely
# - The sheet's used range shrinks from A1:G3 to A1:G2, selection moves to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving record (row 2) with the new token/document number and name.
$ws.Range("A2").Value = 1000588257
$ws.Range("B2").Value = "Stefany Jimenez"
$ws.Range("E2").Value = 1000588257
# F2's CONCAT formula stays in place and will recalc to 1000588257@qmas.com automatically.

# Remove the second record (row 3) completely.
$ws.Rows(3).Delete()

# Drop the (no-op, fill-only) formatting that had been applied to F2.
$ws.Range("F2").ClearFormats()

# Leave the selection on A2, matching the saved workbook view.
$ws.Range("A2").Select()
